$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 250
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H64").Value = 9075.6
$ws.Range("I64").Value = 5005.6665
$ws.Range("K64").Value = 5005.6665
$ws.Range("M64").Value = -4757.6665
$ws.Range("H67").Value = 9075.6
$ws.Range("I67").Value = 5005.6665
$ws.Range("K67").Value = 5005.6665
$ws.Range("M67").Value = -4147.6665
$ws.Range("H74").Value = 9539.362999999999
$ws.Range("I74").Value = 5489.6665
$ws.Range("K74").Value = 5489.6665
$ws.Range("M74").Value = -4553.6665
$ws.Range("H77").Value = 9539.362999999999
$ws.Range("I77").Value = 5489.6665
$ws.Range("K77").Value = 27448.3325
$ws.Range("M77").Value = -22768.3325
$ws.Range("H100").Value = 4385.4
$ws.Range("I100").Value = 4385.4
$ws.Range("K100").Value = 4385.4
$ws.Range("M100").Value = -3844.4
$ws.Range("H131").Value = 73613.2
$ws.Range("I131").Value = 77442.71000000001
$ws.Range("K131").Value = 232328.13
$ws.Range("M131").Value = -227288.13
$ws.Range("H132").Value = 2101
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
$ws.Range("H138").Value = 2721.1292
$ws.Range("J138").Value = 3755.1428
$ws.Range("L138").Value = 11265.4284
$ws.Range("N138").Value = -21545.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 376.63635
$ws.Range("I5").Value = 340.8889
$ws.Range("K5").Value = 340.8889
$ws.Range("M5").Value = -228.8889
$ws.Range("H61").Value = 5383511.5
$ws.Range("I61").Value = 6416110
$ws.Range("K61").Value = 6416110
$ws.Range("M61").Value = -6415898
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("N104").Value = 0
$ws.Range("H132").Value = 5586.115
$ws.Range("I132").Value = 4773.476
$ws.Range("K132").Value = 14320.428
$ws.Range("M132").Value = -11790.428
$ws.Range("H136").Value = 5383511.5
$ws.Range("I136").Value = 6416110
$ws.Range("K136").Value = 19248330
$ws.Range("M136").Value = -19245780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 376.63635
$ws.Range("I4").Value = 340.8889
$ws.Range("K4").Value = 340.8889
$ws.Range("M4").Value = -225.8889
$ws.Range("H86").Value = 8699653
$ws.Range("I86").Value = 3776.1538
$ws.Range("J86").Value = 20004292
$ws.Range("K86").Value = 3776.1538
$ws.Range("L86").Value = 20004292
$ws.Range("M86").Value = -2653.1538
$ws.Range("N86").Value = -20006538
$ws.Range("H89").Value = 8699653
$ws.Range("I89").Value = 3776.1538
$ws.Range("J89").Value = 20004292
$ws.Range("K89").Value = 18880.769
$ws.Range("L89").Value = 100021460
$ws.Range("M89").Value = -13264.769
$ws.Range("N89").Value = -100032692
$ws.Range("H96").Value = 7710
$ws.Range("I96").Value = 7710
$ws.Range("K96").Value = 7710
$ws.Range("M96").Value = -4964
$ws.Range("H107").Value = 2532.16
$ws.Range("I107").Value = 2224.1428
$ws.Range("K107").Value = 2224.1428
$ws.Range("M107").Value = -304.1428000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 5000
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5480
$ws.Range("H27").Value = 5000
$ws.Range("J27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5384
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").ClearContents()
$ws.Range("N51").Value = 0
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("N61").Value = 0
$ws.Range("H62").Value = 8750
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 8750
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 50000
$ws.Range("N65").Value = -56240
$ws.Range("H74").Value = 34700.668
$ws.Range("J74").Value = 34700.668
$ws.Range("L74").Value = 34700.668
$ws.Range("N74").Value = -36448.668
$ws.Range("H77").Value = 34700.668
$ws.Range("J77").Value = 34700.668
$ws.Range("L77").Value = 104102.004
$ws.Range("N77").Value = -112838.004
$ws.Range("H105").Value = 2053.7334
$ws.Range("I105").Value = 1545.8
$ws.Range("K105").Value = 1545.8
$ws.Range("M105").Value = 201.2
$ws.Range("H132").Value = 1866.8182
$ws.Range("I132").Value = 1159.5555
$ws.Range("K132").Value = 3478.6665
$ws.Range("M132").Value = -948.6664999999998
$ws.Range("H134").Value = 4661.8716
$ws.Range("I134").Value = 2004.3704
$ws.Range("K134").Value = 6013.1112
$ws.Range("M134").Value = -3478.1112
$ws.Range("H141").Value = 39184.777
$ws.Range("J141").Value = 39184.777
$ws.Range("L141").Value = 39184.777
$ws.Range("N141").Value = -49544.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1200.4
$ws.Range("I34").Value = 299.5
$ws.Range("J34").Value = 1801
$ws.Range("K34").Value = 898.5
$ws.Range("L34").Value = 5403
$ws.Range("M34").Value = -814.5
$ws.Range("N34").Value = -5571
$ws.Range("H39").Value = 699.8
$ws.Range("J39").Value = 511
$ws.Range("L39").Value = 1533
$ws.Range("N39").Value = -2121
$ws.Range("H55").Value = 1335
$ws.Range("J55").Value = 1627.75
$ws.Range("L55").Value = 4883.25
$ws.Range("N55").Value = -5237.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2945.5715
$ws.Range("I102").Value = 2374
$ws.Range("K102").Value = 2374
$ws.Range("M102").Value = -752
$ws.Range("H113").Value = 801941.2
$ws.Range("I113").Value = 1002176.5
$ws.Range("K113").Value = 1002176.5
$ws.Range("M113").Value = -1000006.5
$ws.Range("H122").Value = 3100.182
$ws.Range("J122").Value = 6599.3335
$ws.Range("L122").Value = 19798.0005
$ws.Range("N122").Value = -24698.0005
$ws.Range("H126").Value = 4908
$ws.Range("I126").Value = 4908
$ws.Range("K126").Value = 14724
$ws.Range("M126").Value = -12254
$ws.Range("H132").Value = 2065.3845
$ws.Range("I132").Value = 2265.375
$ws.Range("J132").Value = 1745.4
$ws.Range("K132").Value = 6796.125
$ws.Range("L132").Value = 5236.200000000001
$ws.Range("M132").Value = -4266.125
$ws.Range("N132").Value = -10296.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4300.8945
$ws.Range("I40").Value = 4357.125
$ws.Range("K40").Value = 4357.125
$ws.Range("M40").Value = -4221.125
$ws.Range("H68").Value = 6330.8335
$ws.Range("I68").Value = 7496.25
$ws.Range("K68").Value = 7496.25
$ws.Range("M68").Value = -6747.25
$ws.Range("H71").Value = 6330.8335
$ws.Range("I71").Value = 7496.25
$ws.Range("K71").Value = 37481.25
$ws.Range("M71").Value = -33737.25
